# Update the weekly class-schedule grid on "Cola aqui os valores" so it
# matches the new best solution found by the genetic-algorithm heuristic.
# Only the subject labels in the two 5-day blocks (B:F and I:M) for the
# class-period rows actually change; everything else (headers, "LANCHE"
# break rows, formulas on the other sheet, etc.) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cola aqui os valores")

# Row 3
$ws.Range("B3").Value = "Geografia"
$ws.Range("C3").Value = "Língua Portuguesa"
$ws.Range("D3").Value = "Língua Portuguesa"
$ws.Range("E3").Value = "Mind Makers"
$ws.Range("F3").Value = "História"
$ws.Range("I3").Value = "Inglês"
$ws.Range("J3").Value = "História"
$ws.Range("L3").Value = "Língua Portuguesa"
$ws.Range("M3").Value = "Música"

# Row 4
$ws.Range("C4").Value = "Língua Portuguesa"
$ws.Range("D4").Value = "Língua Portuguesa"
$ws.Range("E4").Value = "Geografia"
$ws.Range("F4").Value = "História"
$ws.Range("I4").Value = "Língua Portuguesa"
$ws.Range("J4").Value = "Artes"
$ws.Range("L4").Value = "Ed. Financeira"
$ws.Range("M4").Value = "Inglês"

# Row 5
$ws.Range("D5").Value = "Língua Portuguesa"
$ws.Range("E5").Value = "Ed. Financeira"
$ws.Range("F5").Value = "Inglês"
$ws.Range("I5").Value = "Língua Portuguesa"
$ws.Range("J5").Value = "Artes"
$ws.Range("K5").Value = "Matemática"
$ws.Range("L5").Value = "Inglês"
$ws.Range("M5").Value = "Educação Física"

# Row 7
$ws.Range("B7").Value = "Língua Portuguesa"
$ws.Range("C7").Value = "Artes"
$ws.Range("D7").Value = "Ciências"
$ws.Range("E7").Value = "Matemática"
$ws.Range("F7").Value = "Matemática"
$ws.Range("J7").Value = "Matemática"
$ws.Range("K7").Value = "Ensino Religioso"
$ws.Range("L7").Value = "Língua Portuguesa"
$ws.Range("M7").Value = "Língua Portuguesa"

# Row 8
$ws.Range("B8").Value = "Língua Portuguesa"
$ws.Range("C8").Value = "Artes"
$ws.Range("D8").Value = "Ciências"
$ws.Range("E8").Value = "Ensino Religioso"
$ws.Range("F8").Value = "Matemática"
$ws.Range("J8").Value = "Matemática"
$ws.Range("K8").Value = "Mind Makers"
$ws.Range("L8").Value = "Língua Portuguesa"
$ws.Range("M8").Value = "Língua Portuguesa"

# Row 13
$ws.Range("B13").Value = "Língua Portuguesa"
$ws.Range("C13").Value = "Ed. Financeira"
$ws.Range("E13").Value = "Inglês"
$ws.Range("F13").Value = "Inglês"
$ws.Range("I13").Value = "Matemática"
$ws.Range("J13").Value = "Inglês"
$ws.Range("K13").Value = "Mind Makers"
$ws.Range("L13").Value = "Geografia"
$ws.Range("M13").Value = "Matemática"

# Row 14
$ws.Range("B14").Value = "Artes"
$ws.Range("D14").Value = "Matemática"
$ws.Range("E14").Value = "Língua Portuguesa"
$ws.Range("F14").Value = "Língua Portuguesa"
$ws.Range("I14").Value = "Matemática"
$ws.Range("K14").Value = "Geografia"
$ws.Range("L14").Value = "Inglês"

# Row 15
$ws.Range("B15").Value = "Artes"
$ws.Range("C15").Value = "Matemática"
$ws.Range("D15").Value = "História"
$ws.Range("E15").Value = "Geografia"
$ws.Range("F15").Value = "Língua Portuguesa"
$ws.Range("J15").Value = "História"
$ws.Range("K15").Value = "Inglês"
$ws.Range("L15").Value = "Língua Portuguesa"

# Row 17
$ws.Range("C17").Value = "Ciências"
$ws.Range("D17").Value = "Língua Portuguesa"
$ws.Range("E17").Value = "Ensino Religioso"
$ws.Range("F17").Value = "Educação Física"
$ws.Range("I17").Value = "Artes"
$ws.Range("J17").Value = "Língua Portuguesa"
$ws.Range("K17").Value = "Ed. Financeira"
$ws.Range("L17").Value = "Ciências"
$ws.Range("M17").Value = "Música"

# Row 18
$ws.Range("C18").Value = "Ciências"
$ws.Range("E18").Value = "Mind Makers"
$ws.Range("F18").Value = "Música"
$ws.Range("I18").Value = "Artes"
$ws.Range("J18").Value = "Língua Portuguesa"
$ws.Range("K18").Value = "Ensino Religioso"
$ws.Range("L18").Value = "Ciências"
$ws.Range("M18").Value = "Educação Física"
